$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price column (D) values are plain numeric-looking strings in the source data
# (e.g. "1.000", "0.7310") that must stay literal text, matching the original
# inlineStr cells. Force text format before assigning so Excel does not coerce
# them into numbers and silently drop significant trailing zeros.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.797.22"
$ws.Range("E2").Value = "  +4.18%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.871.62"
$ws.Range("E3").Value = "  +3.00%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "277.41"
$ws.Range("E5").Value = "  -0.14%  "
$ws.Range("E6").Value = "  +0.02%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5294"
$ws.Range("E7").Value = "  +3.85%  "
$ws.Range("E8").Value = "  -3.02%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06939"
$ws.Range("E9").Value = "  +3.95%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.03"
$ws.Range("E10").Value = "  -0.13%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.8046"
$ws.Range("E11").Value = "  -2.68%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07764"
$ws.Range("E12").Value = "  -1.76%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.846.83"
$ws.Range("E13").Value = "  +1.65%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "90.33"
$ws.Range("E14").Value = "  +3.03%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.181"
$ws.Range("E15").Value = "  +2.12%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.56"
$ws.Range("E16").Value = "  +3.40%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.000"
$ws.Range("E17").Value = "  +0.03%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008038"
$ws.Range("E18").Value = "  +0.08%  "
$ws.Range("E19").Value = "  -0.02%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "26.826.15"
$ws.Range("E20").Value = "  +4.10%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.102.11"
$ws.Range("E21").Value = "  +2.47%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.03"
$ws.Range("E23").Value = "  +0.37%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.168"
$ws.Range("E24").Value = "  +1.02%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.379"
$ws.Range("E25").Value = "  +6.89%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "146.55"
$ws.Range("E26").Value = "  +3.45%  "
$ws.Range("E27").Value = "  +1.43%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.658"
$ws.Range("E28").Value = "  -0.91%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "113.42"
$ws.Range("E29").Value = "  +3.75%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.346"
$ws.Range("E30").Value = "  +0.43%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.328"
$ws.Range("E31").Value = "  +2.17%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.08911"
$ws.Range("E32").Value = "  +1.60%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04937"
$ws.Range("E33").Value = "  +0.94%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.168"
$ws.Range("E34").Value = "  +2.84%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7310"
$ws.Range("E35").Value = "  +0.44%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.883"
$ws.Range("E36").Value = "  +0.75%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.274"
$ws.Range("E37").Value = "  +4.37%  "
$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01857"
$ws.Range("E38").Value = "  +0.14%  "
$ws.Range("B39").Value = "RenderToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.322"
$ws.Range("E39").Value = "  -2.54%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5147"
$ws.Range("E40").Value = "  -0.34%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9504"
$ws.Range("E41").Value = "  -1.61%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "116.33"
$ws.Range("E42").Value = "  +4.81%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.166"
$ws.Range("E43").Value = "  -0.87%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.113"
$ws.Range("E44").Value = "  +1.08%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.000"
$ws.Range("E45").Value = "  +0.01%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4469"
$ws.Range("E46").Value = "  -2.30%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.1341"
$ws.Range("E47").Value = "  -1.68%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.331"
$ws.Range("E48").Value = "  +1.43%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "36.35"
$ws.Range("E49").Value = "  -0.37%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05946"
$ws.Range("E50").Value = "  +1.82%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.490"
$ws.Range("E51").Value = "  -0.71%  "
